$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# New TPM-based results: update existing rows 2-6 and append new rows 7-11
# for the additional "MuSCs" sending-cluster block. Also renames the old
# "Neutrophils" target cluster to "Resolving-Mac" (and reuses it for the new
# block), and introduces a "MuSCs" sending cluster alongside "FAPs".
# ---------------------------------------------------------------------------

$rows = @(
    # r, A(sending), B(ligand), C(receptor), D(target), E, F, G, H, I, J, K, L, M, N, O, P, Q, R, S, T
    @(2,  "FAPs",  "Wnt5b", "Fzd8", "ECs",               1, 0.3333333333333333, 0.3693360000000001, 1.108008, 0.9453581798061689, 0.9453581798061688, 3, 1,                  2.157506,           6.472517999999999, 0.3549648016839517, 0.3549648016839516, 0.796844636016,       7.171601724144,     0.3355688788151983, 0.3355688788151981),
    @(3,  "FAPs",  "Wnt5b", "Fzd8", "FAPs",              1, 0.3333333333333333, 0.3693360000000001, 1.108008, 0.9453581798061689, 0.9453581798061688, 3, 1,                  2.913654666666667, 8.740964,           0.4793705560628122, 0.4793705560628121, 1.076117559968,       9.685058039712001, 0.4531768763322112, 0.453176876332211),
    @(4,  "FAPs",  "Wnt5b", "Fzd8", "Inflammatory-Mac",  1, 0.3333333333333333, 0.3693360000000001, 1.108008, 0.9453581798061689, 0.9453581798061688, 2, 0.6666666666666666, 0.018986,           0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.007012213296000001, 0.06310991966399999, 0.002952997921296791, 0.00295299792129679),
    @(5,  "FAPs",  "Wnt5b", "Fzd8", "MuSCs",             1, 0.3333333333333333, 0.3693360000000001, 1.108008, 0.9453581798061689, 0.9453581798061688, 3, 1,                  0.9848966666666668, 2.95469,            0.1620406385718132, 0.1620406385718132, 0.3637577952800001,  3.273820157520001, 0.1531864431348786, 0.1531864431348786),
    @(6,  "FAPs",  "Wnt5b", "Fzd8", "Resolving-Mac",     1, 0.3333333333333333, 0.3693360000000001, 1.108008, 0.9453581798061689, 0.9453581798061688, 1, 0.3333333333333333, 0.003041,           0.009122999999999999, 0.0005003221135518961, 0.000500322113551896, 0.001123150776,      0.010108356984,     0.0004729836025841958, 0.0004729836025841957),
    @(7,  "MuSCs", "Wnt5b", "Fzd8", "ECs",               1, 0.3333333333333333, 0.02134766666666667, 0.064043, 0.05464182019383115, 0.05464182019383115, 3, 1,                 2.157506,           6.472517999999999, 0.3549648016839517, 0.3549648016839516, 0.04605771891933333, 0.4145194702739999, 0.01939592286875342, 0.01939592286875342),
    @(8,  "MuSCs", "Wnt5b", "Fzd8", "FAPs",              1, 0.3333333333333333, 0.02134766666666667, 0.064043, 0.05464182019383115, 0.05464182019383115, 3, 1,                 2.913654666666667, 8.740964,           0.4793705560628122, 0.4793705560628121, 0.06219972860577778, 0.559797557452,     0.02619367973060104, 0.02619367973060104),
    @(9,  "MuSCs", "Wnt5b", "Fzd8", "Inflammatory-Mac",  1, 0.3333333333333333, 0.02134766666666667, 0.064043, 0.05464182019383115, 0.05464182019383115, 2, 0.6666666666666666, 0.018986,          0.05695799999999999, 0.003123681567871193, 0.003123681567871192, 0.0004053067993333333, 0.003647761194,   0.0001706836465744023, 0.0001706836465744023),
    @(10, "MuSCs", "Wnt5b", "Fzd8", "MuSCs",             1, 0.3333333333333333, 0.02134766666666667, 0.064043, 0.05464182019383115, 0.05464182019383115, 3, 1,                 0.9848966666666668, 2.95469,            0.1620406385718132, 0.1620406385718132, 0.02102524574111111, 0.18922721167,      0.008854195436934598, 0.008854195436934598),
    @(11, "MuSCs", "Wnt5b", "Fzd8", "Resolving-Mac",     1, 0.3333333333333333, 0.02134766666666667, 0.064043, 0.05464182019383115, 0.05464182019383115, 1, 0.3333333333333333, 0.003041,          0.009122999999999999, 0.0005003221135518961, 0.000500322113551896, 0.00006491825433333333, 0.000584264289, 0.00002733851096770028, 0.00002733851096770027)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value  = $row[1]
    $ws.Cells.Item($r, 2).Value  = $row[2]
    $ws.Cells.Item($r, 3).Value  = $row[3]
    $ws.Cells.Item($r, 4).Value  = $row[4]
    $ws.Cells.Item($r, 5).Value  = $row[5]
    $ws.Cells.Item($r, 6).Value  = $row[6]
    $ws.Cells.Item($r, 7).Value  = $row[7]
    $ws.Cells.Item($r, 8).Value  = $row[8]
    $ws.Cells.Item($r, 9).Value  = $row[9]
    $ws.Cells.Item($r, 10).Value = $row[10]
    $ws.Cells.Item($r, 11).Value = $row[11]
    $ws.Cells.Item($r, 12).Value = $row[12]
    $ws.Cells.Item($r, 13).Value = $row[13]
    $ws.Cells.Item($r, 14).Value = $row[14]
    $ws.Cells.Item($r, 15).Value = $row[15]
    $ws.Cells.Item($r, 16).Value = $row[16]
    $ws.Cells.Item($r, 17).Value = $row[17]
    $ws.Cells.Item($r, 18).Value = $row[18]
    $ws.Cells.Item($r, 19).Value = $row[19]
    $ws.Cells.Item($r, 20).Value = $row[20]
}
